$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: the naive component forecaster had populated spurious one-period-ahead
# y_0_forecast/y_1_forecast values in the very first forecast row(s), and carried
# slightly imprecise values elsewhere. Remove the erroneous cells and correct the
# forecast values to their recomputed precision.

# Remove the incorrectly-populated forecast cells (row 2 had no valid forecast yet,
# and row 3's y_0_forecast was likewise invalid).
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()

# Update remaining forecast values with corrected precision.
$ws.Range("E3").Value = 2.957479223435744

$ws.Range("C4").Value = -0.01587181126745385
$ws.Range("E4").Value = 3.8351443707757

$ws.Range("C5").Value = -0.02256889165886955
$ws.Range("E5").Value = -0.6757980944263275

$ws.Range("C6").Value = 0.09611428386595566

$ws.Range("E7").Value = -0.971238541762387

$ws.Range("C8").Value = -0.001350220946472191
$ws.Range("E8").Value = 0.6008487920565075

$ws.Range("E9").Value = -1.58998093318411

$ws.Range("C10").Value = -0.5761528471665334
$ws.Range("E10").Value = 0.4501721032283301

$ws.Range("E11").Value = -0.150175137493469

$ws.Range("E13").Value = 2.372078088364704

$ws.Range("C14").Value = -0.4278219446121501
$ws.Range("E14").Value = -2.378564786744752

$ws.Range("C15").Value = -1.026566979837429

$ws.Range("C17").Value = 0.4636049209196802

$ws.Range("C18").Value = 0.6216390921348403
$ws.Range("E18").Value = -1.097580983230539

$ws.Range("C19").Value = -0.6768900623516871
